$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new "6th Annual Translational Microbiome Conference" event,
# replacing the old "THE DYNAMIC MICROBIOME" event.
$ws.Range("A2").Value = "MICROBIOME listserv"
$ws.Range("B2").Value = "6th Annual Translational Microbiome Conference"
$ws.Range("C2").Value = Get-Date -Year 2020 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Range("D2").Value = Get-Date -Year 2020 -Month 4 -Day 23 -Hour 0 -Minute 0 -Second 0
$ws.Range("E2").Value = 1587499573
# The new timestamp no longer needs the old integer number format, so drop
# the special formatting back to the plain style used by the cells below it.
$ws.Range("E3").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("F2").Value = "The Westin Copley Place 10 Huntington Avenue Boston, MA 02116 USA"
$ws.Range("G2").Value = "https://microbiomeconference.com"
$ws.Range("H2").Value = "Conference"
$ws.Range("I2").Value = "no"
$ws.Range("J2").Value = "The conference will bring together the leading microbiome companies working to successfully commercialize microbiome-based diagnostics, therapeutics, adjunct therapies and direct-to-consumer services and products across a range of therapeutic areas, from gut to skin and beyond. Expanding on this core focus year on year, the conference continues to offer informative, hands-on workshops, panels sponsored and directed by industry to focus on your concerns and tightly curated content that goes beyond the usual infomercials to provide tangible and useful insights to how your peers are addressing their commercial approach to the space."
$ws.Range("K2").Value = "microbiome,analysis"

# Row 2 is shorter now that the description text is shorter.
$ws.Rows.Item(2).RowHeight = 102

# The trailing blank row (row 7) is no longer needed.
$ws.Rows.Item(7).Delete()

# Move the active selection to A3, matching the author's final cursor position.
$ws.Range("A3").Select()
